$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "67.659.26"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  +3.95%  "
# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.782.88"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  +8.44%  "
# Row 4
$ws.Range("E4").Value = "  +0.32%  "
# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "419.44"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +0.49%  "
# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "132.41"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +1.04%  "
# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.766.27"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  +8.27%  "
# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.647"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  -0.93%  "
# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "1.00"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  +0.01%  "
# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.769"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  -0.63%  "
# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.184"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  +14.52%  "
# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0000412"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  +56.39%  "
# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "42.66"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  -0.64%  "
# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "10.40"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +4.54%  "
# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.372.73"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  +8.34%  "
# Row 16
$ws.Range("E16").Value = "  -0.41%  "
# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.791.46"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  +7.93%  "
# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "20.49"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  +0.93%  "
# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "13.24"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +2.81%  "
# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "1.13"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  +3.10%  "
# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "67.738.43"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  +4.35%  "
# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "445.59"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +0.70%  "
# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "15.82"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +20.55%  "
# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "90.89"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +1.61%  "
# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "3.08"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -4.50%  "
# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "38.17"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  +12.82%  "
# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "3.33"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -0.64%  "
# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.08"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  +1.95%  "
# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.08"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  +4.69%  "
# Row 30
$ws.Range("B30").Value = "Hedera"
$ws.Range("C30").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.124"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  +6.54%  "
# Row 31
$ws.Range("B31").Value = "Cosmos"
$ws.Range("C31").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "12.64"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  +1.23%  "
# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.77"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -1.44%  "
# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "7.17"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  -3.14%  "
# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.164"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  +1.90%  "
# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "41.48"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  +5.44%  "
# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "57.83"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  +1.72%  "
# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.999"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  +0.00%  "
# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0490"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  -2.65%  "
# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0₃0730"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +1.38%  "
# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.00"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +29.62%  "
# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.148"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +0.34%  "
# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "28.30"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  +32.45%  "
# Row 43
$ws.Range("E43").Value = "  +0.19%  "
# Row 44
$ws.Range("E44").Value = "  +4.09%  "
# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "148.31"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +1.23%  "
# Row 46
$ws.Range("B46").Value = "ApeXProtocol"
$ws.Range("C46").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.19"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +25.00%  "
# Row 47
$ws.Range("B47").Value = "ARBITRUM"
$ws.Range("C47").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.12"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  +6.46%  "
# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.89"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -4.10%  "
# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.63"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -5.26%  "
# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "4.30"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -3.96%  "
# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.305"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -1.36%  "
